$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44211
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 2800
$ws.Range("O2").Value = 2800
$ws.Range("P2").Value = 2800
$ws.Range("S2").Value = 1400

# Row 3
$ws.Range("D3").Value = 44211
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 2600
$ws.Range("O3").Value = 2600
$ws.Range("P3").Value = 2600
$ws.Range("R3").Value = 'Provincia de Linares'
$ws.Range("S3").Value = 1300

# Row 4
$ws.Range("D4").Value = 44204
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 3000
$ws.Range("P4").Value = 3000
$ws.Range("S4").Value = 1500

# Row 5
$ws.Range("D5").Value = 44204
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 2400
$ws.Range("O5").Value = 2400
$ws.Range("P5").Value = 2400
$ws.Range("S5").Value = 1200

# Row 6
$ws.Range("D6").Value = 44232

# Row 7
$ws.Range("D7").Value = 44166
$ws.Range("M7").Value = 1500

# Row 8
$ws.Range("D8").Value = 44162
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 4000
$ws.Range("O8").Value = 4000
$ws.Range("P8").Value = 4000
$ws.Range("S8").Value = 2000

# Row 9
$ws.Range("D9").Value = 44202
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("R9").Value = 'Provincia de Linares'
$ws.Range("S9").Value = 1500

# Row 10
$ws.Range("D10").Value = 44202
$ws.Range("L10").Value = 'Segunda'
$ws.Range("M10").Value = 20
$ws.Range("N10").Value = 2600
$ws.Range("O10").Value = 2600
$ws.Range("P10").Value = 2600
$ws.Range("S10").Value = 1300

# Row 11
$ws.Range("D11").Value = 44165
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 400
$ws.Range("N11").Value = 3400
$ws.Range("O11").Value = 3400
$ws.Range("P11").Value = 3400
$ws.Range("R11").Value = 'Región de O''Higgins'
$ws.Range("S11").Value = 1700

# Row 12
$ws.Range("D12").Value = 44169
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 400
$ws.Range("N12").Value = 3600
$ws.Range("O12").Value = 3600
$ws.Range("P12").Value = 3600
$ws.Range("S12").Value = 1800

# Row 13
$ws.Range("D13").Value = 44172
$ws.Range("M13").Value = 300
$ws.Range("N13").Value = 3400
$ws.Range("O13").Value = 3600
$ws.Range("P13").Value = 3467
$ws.Range("Q13").Value = '$/bandeja 2 kilos'
$ws.Range("R13").Value = 'Provincia de Linares'
$ws.Range("S13").Value = 1734
$ws.Range("T13").Value = 2

# Row 14
$ws.Range("D14").Value = 44187
$ws.Range("M14").Value = 110
$ws.Range("N14").Value = 2600
$ws.Range("O14").Value = 3000
$ws.Range("P14").Value = 2782
$ws.Range("S14").Value = 1391

# Row 15
$ws.Range("D15").Value = 44200
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 50
$ws.Range("N15").Value = 2600
$ws.Range("O15").Value = 2600
$ws.Range("P15").Value = 2600
$ws.Range("S15").Value = 1300

# Row 16
$ws.Range("D16").Value = 44265
$ws.Range("M16").Value = 70
$ws.Range("N16").Value = 3600
$ws.Range("O16").Value = 3800
$ws.Range("P16").Value = 3714
$ws.Range("S16").Value = 1857

# Row 17
$ws.Range("D17").Value = 44186
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 3000
$ws.Range("O17").Value = 3000
$ws.Range("P17").Value = 3000
$ws.Range("R17").Value = 'Provincia de Limarí'
$ws.Range("S17").Value = 1500

# Row 18
$ws.Range("D18").Value = 44235
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 3000
$ws.Range("O18").Value = 3000
$ws.Range("P18").Value = 3000
$ws.Range("S18").Value = 1500

# Row 19
$ws.Range("D19").Value = 44264
$ws.Range("M19").Value = 110
$ws.Range("N19").Value = 3500
$ws.Range("O19").Value = 4000
$ws.Range("P19").Value = 3727
$ws.Range("R19").Value = 'Provincia de Linares'
$ws.Range("S19").Value = 1864

# Row 20
$ws.Range("D20").Value = 44210
$ws.Range("L20").Value = 'Segunda'
$ws.Range("M20").Value = 150
$ws.Range("N20").Value = 2700
$ws.Range("O20").Value = 2700
$ws.Range("P20").Value = 2700
$ws.Range("S20").Value = 1350

# Row 21
$ws.Range("D21").Value = 44176
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 150
$ws.Range("N21").Value = 3500
$ws.Range("O21").Value = 3500
$ws.Range("P21").Value = 3500
$ws.Range("Q21").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 2333
$ws.Range("T21").Value = 1.5

# Row 22
$ws.Range("D22").Value = 44167
$ws.Range("M22").Value = 500
$ws.Range("N22").Value = 3600
$ws.Range("O22").Value = 3600
$ws.Range("P22").Value = 3600
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1800
